# Weekly update: a new week's price record is published for
# Agrícola del Norte S.A. de Arica - Acelga.
# This pushes every existing record down by one row (row 7 becomes row 8,
# row 8 becomes row 9, ... row 67 becomes row 68) and the brand-new
# week's data is written into the now-empty row 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 - Excel shifts rows 7:67 down to 8:68
# and extends the used range to A1:R68, carrying the row's formatting
# (e.g. the date-style D column) down with it.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with this week's record.
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = "7/19/2022"
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 100112009
$ws.Range("G7").Value = "Acelga"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Segunda"
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 1200
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = 1350
$ws.Range("N7").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O7").Value = "Región de Arica y Parinacota"
$ws.Range("P7").Value = 450
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = "Hortaliza"
